$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 3 with login credentials, matching style of existing "quotePrefix" cells (style index 1)
$ws.Range("A3").Value = "super@admin.com"
$ws.Range("B3").Value = "password"

# Copy style from existing cells that use style index 1 (e.g. A5) to the new cells
[void]$ws.Range("A5").Copy()
[void]$ws.Range("A3").PasteSpecial(-4122) # xlPasteFormats
[void]$ws.Range("A5").Copy()
[void]$ws.Range("B3").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = $false | Out-Null

# Update the active selection to B3
[void]$ws.Range("B3").Select()
